$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44461
$ws.Range("L2").Value = "Especial"
$ws.Range("N2").Value = 31000
$ws.Range("O2").Value = 32000
$ws.Range("P2").Value = 31500
$ws.Range("S2").Value = 3150

$ws.Range("D3").Value = 44461
$ws.Range("M3").Value = 30
$ws.Range("N3").Value = 30000
$ws.Range("O3").Value = 30000
$ws.Range("P3").Value = 30000
$ws.Range("S3").Value = 3000

$ws.Range("D4").Value = 44446

$ws.Range("D7").Value = 44448

$ws.Range("D8").Value = 44452
$ws.Range("L8").Value = "Primera"
$ws.Range("N8").Value = 21000
$ws.Range("O8").Value = 22000
$ws.Range("P8").Value = 21500
$ws.Range("S8").Value = 2150

$ws.Range("D9").Value = 44447
$ws.Range("M9").Value = 60
$ws.Range("N9").Value = 21000
$ws.Range("O9").Value = 22000
$ws.Range("P9").Value = 21500
$ws.Range("S9").Value = 2150
